# Maestro.xlsx edit: insert a new "Cif desengrasante" detergent row right
# after the header (becomes new row 2, shifting every existing product row
# down by one), and append a brand-new "Campanita" papel higienico row at
# the end of the table (new row 34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a fresh row at row 2 (shifts rows 2..32 down to 3..33).
# ---------------------------------------------------------------------
$ws.Range("A2:O2").Insert()

# Fill the new row 2 with the new product's data. Columns are written in
# left-to-right order so any brand-new shared strings are appended to
# sharedStrings.xml in the same order the target workbook uses.
$ws.Range("A2").Value2 = 7791290794115
$ws.Range("B2").Value2 = "Detergente"
$ws.Range("C2").Value2 = "desengrasante"
$ws.Range("D2").Value2 = "bioactive lima"
$ws.Range("E2").Value2 = "Cif"
$ws.Range("F2").Value2 = 500
$ws.Range("G2").Value2 = "ml."
$ws.Range("H2").Value2 = "Botella"
$ws.Range("I2").Value2 = "Detergentes"
$ws.Range("J2").Value2 = "Argentina"
$ws.Range("K2").Value2 = 6
$ws.Range("L2").Value2 = $false
$ws.Range("M2").Value2 = $true
$ws.Range("N2").Value2 = "C:\VentaSoft\Imágenes de artículos\7791290794115.png"
$ws.Range("O2").Value2 = $true

# The freshly inserted row picked up a generic "inserted row" format; put
# the usual bordered/right-aligned "ImagenExactaDelArticulo" format (the
# one every other data row in O2:O25 uses) back onto O2 by copying it from
# O3 (the row right below, which still carries the original formatting).
$ws.Range("O3").Copy()
$ws.Range("O2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Append a brand-new row after the last existing product row (old row
#    32 is now row 33, so the new row goes at row 34).
# ---------------------------------------------------------------------
# Give A34 the same integer-style formatting ("s=1") used by the rest of
# column A by copying the format from A33 before filling in the value.
$ws.Range("A33").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A34").Value2 = 7791070000078
$ws.Range("B34").Value2 = "Papel higiénico"
$ws.Range("C34").Value2 = "simple hoja"
$ws.Range("D34").Value2 = '"soft"'
$ws.Range("E34").Value2 = "Campanita"
$ws.Range("F34").Value2 = 4
$ws.Range("G34").Value2 = "und."
$ws.Range("H34").Value2 = "Bolsa"
$ws.Range("I34").Value2 = "Papeles Higiénicos"
$ws.Range("J34").Value2 = "Argentina"
$ws.Range("K34").Value2 = 10
$ws.Range("L34").Value2 = $false
$ws.Range("M34").Value2 = $false
$ws.Range("N34").Value2 = "C:\VentaSoft\Imágenes de artículos\7791070000078.png"
$ws.Range("O34").Value2 = $true
